$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1620.1666
$ws.Cells.Item(15, 9).Value = 1620.1666
$ws.Cells.Item(15, 11).Value = 4860.4998
$ws.Cells.Item(15, 13).Value = -4691.4998

$ws.Cells.Item(88, 8).Value = 1909
$ws.Cells.Item(88, 9).Value = 966
$ws.Cells.Item(88, 10).Value = 2262.625
$ws.Cells.Item(88, 11).Value = 966
$ws.Cells.Item(88, 12).Value = 2262.625
$ws.Cells.Item(88, 13).Value = -560
$ws.Cells.Item(88, 14).Value = -3074.625

$ws.Cells.Item(91, 8).Value = 1909
$ws.Cells.Item(91, 9).Value = 966
$ws.Cells.Item(91, 10).Value = 2262.625
$ws.Cells.Item(91, 11).Value = 966
$ws.Cells.Item(91, 12).Value = 2262.625
$ws.Cells.Item(91, 13).Value = 438
$ws.Cells.Item(91, 14).Value = -5070.625

$ws.Cells.Item(100, 8).Value = 2443.889
$ws.Cells.Item(100, 9).Value = 2349.1667
$ws.Cells.Item(100, 10).Value = 2633.3333
$ws.Cells.Item(100, 11).Value = 2349.1667
$ws.Cells.Item(100, 12).Value = 2633.3333
$ws.Cells.Item(100, 13).Value = -1808.1667
$ws.Cells.Item(100, 14).Value = -3715.3333

$ws.Cells.Item(137, 8).Value = 2993.6667
$ws.Cells.Item(137, 9).Value = 1475.5
$ws.Cells.Item(137, 11).Value = 4426.5
$ws.Cells.Item(137, 13).Value = -1876.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 100002584
$ws.Cells.Item(102, 9).Value = 2868.889
$ws.Cells.Item(102, 10).Value = 1000000000
$ws.Cells.Item(102, 11).Value = 2868.889
$ws.Cells.Item(102, 12).Value = 1000000000
$ws.Cells.Item(102, 13).Value = -1246.889
$ws.Cells.Item(102, 14).Value = -1000003244

$ws.Cells.Item(134, 8).Value = 49800
$ws.Cells.Item(134, 10).Value = 49800
$ws.Cells.Item(134, 12).Value = 49800
$ws.Cells.Item(134, 14).Value = -59940

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(56, 8).Value = 13527.5
$ws.Cells.Item(56, 10).Value = 13527.5
$ws.Cells.Item(56, 12).Value = 13527.5
$ws.Cells.Item(56, 14).Value = -15005.5

$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).ClearContents()

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).ClearContents()

$ws.Cells.Item(86, 8).Value = 2059.348
$ws.Cells.Item(86, 9).Value = 2306.4
$ws.Cells.Item(86, 10).Value = 1596.125
$ws.Cells.Item(86, 11).Value = 2306.4
$ws.Cells.Item(86, 12).Value = 1596.125
$ws.Cells.Item(86, 13).Value = -1183.4
$ws.Cells.Item(86, 14).Value = -3842.125

$ws.Cells.Item(89, 8).Value = 2059.348
$ws.Cells.Item(89, 9).Value = 2306.4
$ws.Cells.Item(89, 10).Value = 1596.125
$ws.Cells.Item(89, 11).Value = 11532
$ws.Cells.Item(89, 12).Value = 7980.625
$ws.Cells.Item(89, 13).Value = -5916
$ws.Cells.Item(89, 14).Value = -19212.625

$ws.Cells.Item(99, 8).Value = 1984.4445
$ws.Cells.Item(99, 9).Value = 1149.8
$ws.Cells.Item(99, 10).Value = 3027.75
$ws.Cells.Item(99, 11).Value = 1149.8
$ws.Cells.Item(99, 12).Value = 3027.75
$ws.Cells.Item(99, 13).Value = 348.2
$ws.Cells.Item(99, 14).Value = -6023.75

$ws.Cells.Item(105, 8).Value = 3958.7917
$ws.Cells.Item(105, 9).Value = 2357.7856
$ws.Cells.Item(105, 10).Value = 6200.2
$ws.Cells.Item(105, 11).Value = 2357.7856
$ws.Cells.Item(105, 12).Value = 6200.2
$ws.Cells.Item(105, 13).Value = -610.7856000000002
$ws.Cells.Item(105, 14).Value = -9694.200000000001

$ws.Cells.Item(134, 8).Value = 2704.2903
$ws.Cells.Item(134, 9).Value = 2508.3809
$ws.Cells.Item(134, 10).Value = 3115.7
$ws.Cells.Item(134, 11).Value = 7525.1427
$ws.Cells.Item(134, 12).Value = 9347.099999999999
$ws.Cells.Item(134, 13).Value = -4990.1427
$ws.Cells.Item(134, 14).Value = -14417.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 587.75
$ws.Cells.Item(22, 9).Value = 650.3333
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 650.3333
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = -300.3333
$ws.Cells.Item(22, 14).Value = -1100

$ws.Cells.Item(31, 8).Value = 29513.475
$ws.Cells.Item(31, 9).Value = 1775.1428
$ws.Cells.Item(31, 10).Value = 107180.8
$ws.Cells.Item(31, 11).Value = 1775.1428
$ws.Cells.Item(31, 12).Value = 107180.8
$ws.Cells.Item(31, 13).Value = -1480.1428
$ws.Cells.Item(31, 14).Value = -107770.8

$ws.Cells.Item(34, 8).Value = 29513.475
$ws.Cells.Item(34, 9).Value = 1775.1428
$ws.Cells.Item(34, 10).Value = 107180.8
$ws.Cells.Item(34, 11).Value = 1775.1428
$ws.Cells.Item(34, 12).Value = 107180.8
$ws.Cells.Item(34, 13).Value = -1573.1428
$ws.Cells.Item(34, 14).Value = -107584.8

$ws.Cells.Item(62, 8).Value = 2833.3333
$ws.Cells.Item(62, 9).Value = 3000
$ws.Cells.Item(62, 10).Value = 2500
$ws.Cells.Item(62, 11).Value = 3000
$ws.Cells.Item(62, 12).Value = 2500
$ws.Cells.Item(62, 13).Value = -2376
$ws.Cells.Item(62, 14).Value = -3748

$ws.Cells.Item(65, 8).Value = 2833.3333
$ws.Cells.Item(65, 9).Value = 3000
$ws.Cells.Item(65, 10).Value = 2500
$ws.Cells.Item(65, 11).Value = 15000
$ws.Cells.Item(65, 12).Value = 12500
$ws.Cells.Item(65, 13).Value = -11880
$ws.Cells.Item(65, 14).Value = -18740

$ws.Cells.Item(81, 8).Value = 30333.334
$ws.Cells.Item(81, 10).Value = 30333.334
$ws.Cells.Item(81, 12).Value = 30333.334
$ws.Cells.Item(81, 14).Value = -32329.334

$ws.Cells.Item(82, 8).Value = 26900
$ws.Cells.Item(82, 9).Value = 9800
$ws.Cells.Item(82, 11).Value = 9800
$ws.Cells.Item(82, 13).Value = -9439

$ws.Cells.Item(84, 8).Value = 30333.334
$ws.Cells.Item(84, 10).Value = 30333.334
$ws.Cells.Item(84, 12).Value = 91000.00199999999
$ws.Cells.Item(84, 14).Value = -100984.002

$ws.Cells.Item(85, 8).Value = 26900
$ws.Cells.Item(85, 9).Value = 9800
$ws.Cells.Item(85, 11).Value = 9800
$ws.Cells.Item(85, 13).Value = -8552

$ws.Cells.Item(87, 8).Value = 25866.666
$ws.Cells.Item(87, 10).Value = 25866.666
$ws.Cells.Item(87, 12).Value = 25866.666
$ws.Cells.Item(87, 14).Value = -28238.666

$ws.Cells.Item(90, 8).Value = 25866.666
$ws.Cells.Item(90, 10).Value = 25866.666
$ws.Cells.Item(90, 12).Value = 77599.99800000001
$ws.Cells.Item(90, 14).Value = -89455.99800000001

$ws.Cells.Item(105, 8).Value = 4286.6665
$ws.Cells.Item(105, 9).Value = 1775
$ws.Cells.Item(105, 10).Value = 6296
$ws.Cells.Item(105, 11).Value = 1775
$ws.Cells.Item(105, 12).Value = 6296
$ws.Cells.Item(105, 13).Value = -28
$ws.Cells.Item(105, 14).Value = -9790

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 10094.615
$ws.Cells.Item(64, 9).Value = 806
$ws.Cells.Item(64, 10).Value = 11783.454
$ws.Cells.Item(64, 11).Value = 2418
$ws.Cells.Item(64, 12).Value = 35350.362
$ws.Cells.Item(64, 13).Value = -2148
$ws.Cells.Item(64, 14).Value = -35890.362

$ws.Cells.Item(67, 8).Value = 10094.615
$ws.Cells.Item(67, 9).Value = 806
$ws.Cells.Item(67, 10).Value = 11783.454
$ws.Cells.Item(67, 11).Value = 2418
$ws.Cells.Item(67, 12).Value = 35350.362
$ws.Cells.Item(67, 13).Value = -1482
$ws.Cells.Item(67, 14).Value = -37222.362

$ws.Cells.Item(70, 8).Value = 6474.4165
$ws.Cells.Item(70, 9).Value = 2583.3333
$ws.Cells.Item(70, 10).Value = 7771.4443
$ws.Cells.Item(70, 11).Value = 7749.999899999999
$ws.Cells.Item(70, 12).Value = 23314.3329
$ws.Cells.Item(70, 13).Value = -7434.999899999999
$ws.Cells.Item(70, 14).Value = -23944.3329

$ws.Cells.Item(73, 8).Value = 6474.4165
$ws.Cells.Item(73, 9).Value = 2583.3333
$ws.Cells.Item(73, 10).Value = 7771.4443
$ws.Cells.Item(73, 11).Value = 7749.999899999999
$ws.Cells.Item(73, 12).Value = 23314.3329
$ws.Cells.Item(73, 13).Value = -6657.999899999999
$ws.Cells.Item(73, 14).Value = -25498.3329

$ws.Cells.Item(76, 8).Value = 8000
$ws.Cells.Item(76, 9).Value = 4000
$ws.Cells.Item(76, 10).Value = 8800
$ws.Cells.Item(76, 11).Value = 12000
$ws.Cells.Item(76, 12).Value = 26400
$ws.Cells.Item(76, 13).Value = -11617
$ws.Cells.Item(76, 14).Value = -27166

$ws.Cells.Item(79, 8).Value = 8000
$ws.Cells.Item(79, 9).Value = 4000
$ws.Cells.Item(79, 10).Value = 8800
$ws.Cells.Item(79, 11).Value = 12000
$ws.Cells.Item(79, 12).Value = 26400
$ws.Cells.Item(79, 13).Value = -10674
$ws.Cells.Item(79, 14).Value = -29052

$ws.Cells.Item(82, 8).Value = 9699.5
$ws.Cells.Item(82, 9).Value = 299
$ws.Cells.Item(82, 10).Value = 12833
$ws.Cells.Item(82, 11).Value = 897
$ws.Cells.Item(82, 12).Value = 38499
$ws.Cells.Item(82, 13).Value = -491
$ws.Cells.Item(82, 14).Value = -39311

$ws.Cells.Item(85, 8).Value = 9699.5
$ws.Cells.Item(85, 9).Value = 299
$ws.Cells.Item(85, 10).Value = 12833
$ws.Cells.Item(85, 11).Value = 897
$ws.Cells.Item(85, 12).Value = 38499
$ws.Cells.Item(85, 13).Value = 507
$ws.Cells.Item(85, 14).Value = -41307

$ws.Cells.Item(86, 8).Value = 792.73334
$ws.Cells.Item(86, 10).Value = 965
$ws.Cells.Item(86, 12).Value = 2895
$ws.Cells.Item(86, 14).Value = -5267

$ws.Cells.Item(88, 8).Value = 6499.8335
$ws.Cells.Item(88, 10).Value = 6499.8335
$ws.Cells.Item(88, 12).Value = 19499.5005
$ws.Cells.Item(88, 14).Value = -20355.5005

$ws.Cells.Item(89, 8).Value = 792.73334
$ws.Cells.Item(89, 10).Value = 965
$ws.Cells.Item(89, 12).Value = 8685
$ws.Cells.Item(89, 14).Value = -20541

$ws.Cells.Item(91, 8).Value = 6499.8335
$ws.Cells.Item(91, 10).Value = 6499.8335
$ws.Cells.Item(91, 12).Value = 19499.5005
$ws.Cells.Item(91, 14).Value = -22463.5005

$ws.Cells.Item(118, 8).Value = 30305656
$ws.Cells.Item(118, 9).Value = 55556020
$ws.Cells.Item(118, 10).Value = 5220.8
$ws.Cells.Item(118, 11).Value = 166668060
$ws.Cells.Item(118, 12).Value = 15662.4
$ws.Cells.Item(118, 13).Value = -166666817
$ws.Cells.Item(118, 14).Value = -18148.4

$ws.Cells.Item(132, 9).Value = 368.18182
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 3313.63638
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = -783.6363799999999
$ws.Cells.Item(132, 14).Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 30367.6
$ws.Cells.Item(70, 9).Value = 53910.19
$ws.Cells.Item(70, 10).Value = 4346.8423
$ws.Cells.Item(70, 11).Value = 53910.19
$ws.Cells.Item(70, 12).Value = 4346.8423
$ws.Cells.Item(70, 13).Value = -53640.19
$ws.Cells.Item(70, 14).Value = -4886.8423

$ws.Cells.Item(73, 8).Value = 30367.6
$ws.Cells.Item(73, 9).Value = 53910.19
$ws.Cells.Item(73, 10).Value = 4346.8423
$ws.Cells.Item(73, 11).Value = 53910.19
$ws.Cells.Item(73, 12).Value = 4346.8423
$ws.Cells.Item(73, 13).Value = -52974.19
$ws.Cells.Item(73, 14).Value = -6218.8423

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2509.9
$ws.Cells.Item(132, 9).Value = 1980
$ws.Cells.Item(132, 10).Value = 4099.6
$ws.Cells.Item(132, 11).Value = 5940
$ws.Cells.Item(132, 12).Value = 12298.8
$ws.Cells.Item(132, 13).Value = -3410
$ws.Cells.Item(132, 14).Value = -17358.8

$ws.Cells.Item(133, 8).Value = 47845
$ws.Cells.Item(133, 10).Value = 47845
$ws.Cells.Item(133, 12).Value = 47845
$ws.Cells.Item(133, 14).Value = -52905
